# Multi sheet importer (parser) v.0.91.
#
# 1. Rename sheet "Basic_action_types" -> "Basic_chain_product_types".
# 2. On Basic_chain_product_types: add an (empty, but styled) B4 cell,
#    and leave the sheet's selection at F22 (no longer the active tab).
# 3. On Basic_categories: fix the typo'd quantity in B10 (10000 -> 1000),
#    make it the active tab, with the selection left at B10.

$wb = $excel.ActiveWorkbook

$chains = $wb.Worksheets.Item("Basic_chains")
$productTypes = $wb.Worksheets.Item("Basic_action_types")
$productTypes.Name = "Basic_chain_product_types"

# B4 was blank in the source data too -- copy the matching blank-cell
# style (Basic_chains!C2) over rather than inventing a new one.
$chains.Range("C2").Copy($productTypes.Range("B4")) | Out-Null

$productTypes.Activate() | Out-Null
$productTypes.Range("F22").Select() | Out-Null

$categories = $wb.Worksheets.Item("Basic_categories")
$categories.Range("B10").Value = 1000

$categories.Activate() | Out-Null
$categories.Range("B10").Select() | Out-Null
